# wallet.docx template update
#   - i18n key separator ":" -> "." for the {{t:...}} placeholders
#   - new FOR-loop block over diagnosisParagraphs appended after the
#     medications loop
#   - Normal style gets explicit hyphenation/spacing/alignment
#   - two new "(user)" paragraph styles are added (mirrors of the
#     existing Ueberschrift / Verzeichnis styles)
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Fix the translation-key placeholders: "{{t:...}}" -> "{{t....}}"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("{{t:notfallpass.title}}", $false, $false, $false, $false, $false, $true, 1, $false, "{{t.notfallpass.title}}", 2) | Out-Null
$d.Content.Find.Execute("{{t:notfallpass.section.contacts.title}}", $false, $false, $false, $false, $false, $true, 1, $false, "{{t.notfallpass.section.contacts.title}}", 2) | Out-Null
$d.Content.Find.Execute("{{t:notfallpass.section.symptoms.title}}", $false, $false, $false, $false, $false, $true, 1, $false, "{{t.notfallpass.section.symptoms.title}}", 2) | Out-Null
$d.Content.Find.Execute("{{t:notfallpass.section.medications.title}}", $false, $false, $false, $false, $false, $true, 1, $false, "{{t.notfallpass.section.medications.title}}", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Append a new "diagnosisParagraphs" loop right after the
#    medications loop body, closing the medications loop first.
# ---------------------------------------------------------------------
# Locate "{{m.name}} | {{m.dosage}} | {{m.schedule}}" - it is the
# paragraph right before the final "{{END-FOR}}".
$medBody = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7) -eq "{{m.name}} | {{m.dosage}} | {{m.schedule}}") {
        $medBody = $d.Paragraphs($i)
    }
}

$anchor = $medBody.Range
$anchor.Collapse(0)

$anchor.InsertParagraphAfter()
$anchor.Collapse(0)
$anchor.MoveStart(1, 1)
$p1 = $d.Paragraphs($medBody.Index + 1)
$p1.Range.Text = "{{END-FOR}}"

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs($medBody.Index + 2)
$p2.Range.Text = "{{FOR p IN diagnosisParagraphs}}"

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs($medBody.Index + 3)
$p3.Range.Text = "{{p}}"

# ---------------------------------------------------------------------
# 3. Normal style: hyphenation off, zero spacing, left alignment
# ---------------------------------------------------------------------
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.Hyphenation = $false
$normal.ParagraphFormat.SpaceBefore = 0
$normal.ParagraphFormat.SpaceAfter = 0
$normal.ParagraphFormat.Alignment = 0

# ---------------------------------------------------------------------
# 4. New paragraph styles "Überschrift (user)" / "Verzeichnis (user)"
#    - exact mirrors of the built-in "Überschrift" / "Verzeichnis"
#      styles, just with a distinct id/name.
# ---------------------------------------------------------------------
$headingUser = $d.Styles.Add("berschriftuser", 1)
$headingUser.NameLocal = "\u00dcberschrift (user)"
$headingUser.BaseStyle = $d.Styles("Normal")
$headingUser.NextParagraphStyle = $d.Styles("BodyText")
$headingUser.QuickStyle = $true
$headingUser.ParagraphFormat.KeepWithNext = $true
$headingUser.ParagraphFormat.SpaceBefore = 12
$headingUser.ParagraphFormat.SpaceAfter = 6
$headingUser.Font.Name = "Liberation Sans"
$headingUser.Font.NameFarEast = "Microsoft YaHei"
$headingUser.Font.NameBi = "Arial"
$headingUser.Font.Size = 14
$headingUser.Font.SizeBi = 14

$indexUser = $d.Styles.Add("Verzeichnisuser", 1)
$indexUser.NameLocal = "Verzeichnis (user)"
$indexUser.BaseStyle = $d.Styles("Normal")
$indexUser.QuickStyle = $true
$indexUser.ParagraphFormat.NoLineNumber = $true
$indexUser.Font.NameBi = "Arial"
